# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# 1) Update the "last updated" timestamp in A1 (17:22 -> 17:52)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 17:52"

# 2) Swap the "La Palma" / "Lanzarote" rows (row 56 and row 57):
#    Ciudad (A) and Muertes (E) values are swapped; B/C/D are identical
#    between the two rows so nothing else needs to change.
$ws.Range("A56").Value = "Lanzarote"
$ws.Range("E56").Value = 3

$ws.Range("A57").Value = "La Palma"
$ws.Range("E57").Value = 4
